$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The paragraph that narrates the Power of Attorney revocation has a
# Jinja/docassemble conditional:
#
#   {% if property_replace_agent == True %} {{ new_property_agent... }}
#
# with a stray space right after the opening tag, and the (hidden)
# _GoBack bookmark is currently sitting much further along, right after
# the *second* "{% if property_replace_agent == False" tag. The May 16
# feedback moves that bookmark up to sit immediately after the first
# tag's closing "%}" and removes the now-unwanted space that used to
# separate the tag from the text that follows it.
# ------------------------------------------------------------------

$openTag = "{% if property_replace_agent == True %}"

$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute($openTag, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the '$openTag' template tag"
}

# $rng now exactly spans the found tag text; collapse to its end point,
# i.e. the position right after the closing "%}".
$tagEnd = $rng.End

# Re-seat the hidden _GoBack bookmark there (Bookmarks.Add with an
# existing name moves the bookmark rather than creating a duplicate).
$target = $d.Range($tagEnd, $tagEnd)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null

# Delete the single space that used to separate "%}" from the
# following "{{ new_property_agent...".
$spaceRng = $d.Range($tagEnd, $tagEnd + 1)
if ($spaceRng.Text -eq " ") {
    $spaceRng.Delete()
}
